# Apply the "ibm.com -> linuxforhealth.org" rebrand edit to the
# StructureDefinition-match-confidence-level workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (key/value property table) ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-level"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements" (element definition table) ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set on the valueCoding slice row (row 7, column Y)
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/match-confidence-level"

# The longer URL text widens the (bestFit) column.
$elements.Columns.Item(25).ColumnWidth = 62.23828125

# The root "Extension" row (row 2) no longer repeats the ele-1/ext-1
# constraint text in the Constraint(s) column; it now only shows on the
# "Extension.extension" row (row 4, which already carries the same text).
$elements.Range("AI2").Value = ""

$wb.Save()
